# Updated cryptos list on Mon Apr 24 06:07:34 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to stay text (column D/E hold numeric-looking strings
# such as "27.910.98" or "1.009" that Excel would otherwise coerce into
# real numbers); rows 18/19 additionally swap their Coin/Link values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.910.98'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.875.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.61'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4759'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.77%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3932'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.16'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.012'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.81'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.871.81'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.023'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.192'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.011'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001049'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06699'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.03'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.891.78'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.501'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.95'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.334'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.101.80'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.26'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.78'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.097'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.455'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.34'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9739'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09512'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.624'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.325'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.353'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06070'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02228'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.205'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.182'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.009'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5957'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1891'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.256'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5644'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.15'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.927'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06781'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.17'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.36%  '
